$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B26").Value = 6481
$ws.Range("C26").Value = 1011
$ws.Range("D26").Value = 6040734
$ws.Range("E26").Value = 932.0681993519519
$ws.Range("F26").Value = 9.53185736014872
$ws.Range("G26").Value = 7.32484076433122
$ws.Range("H26").Value = 25.79848906050568
